$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data was refreshed with a new week of prices. Two new rows of
# data (for 2023-11-09) are inserted right above the existing "Terminal La
# Palmera de La Serena - Naranja" block that starts at row 1060, pushing the
# rest of that block (and everything below it) down by two rows.
$ws.Rows("1060:1061").Insert()

# Row 1060: Naranja / Cara cara / Primera
$ws.Cells.Item(1060, 1).Value2 = 8
$ws.Cells.Item(1060, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1060, 3).Value2 = "Coquimbo"
$ws.Cells.Item(1060, 4).Value2 = 45239
$ws.Cells.Item(1060, 5).Value2 = 4
$ws.Cells.Item(1060, 6).Value2 = "Fruta"
$ws.Cells.Item(1060, 7).Value2 = 100102
$ws.Cells.Item(1060, 8).Value2 = "Cítricos"
$ws.Cells.Item(1060, 9).Value2 = 100102005
$ws.Cells.Item(1060, 10).Value2 = "Naranja"
$ws.Cells.Item(1060, 11).Value2 = "Cara cara"
$ws.Cells.Item(1060, 12).Value2 = "Primera"
$ws.Cells.Item(1060, 13).Value2 = 10
$ws.Cells.Item(1060, 14).Value2 = 180000
$ws.Cells.Item(1060, 15).Value2 = 190000
$ws.Cells.Item(1060, 16).Value2 = 185000
$ws.Cells.Item(1060, 17).Value2 = "`$/bins (400 kilos)"
$ws.Cells.Item(1060, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(1060, 19).Value2 = 462
$ws.Cells.Item(1060, 20).Value2 = 400

# Row 1061: Naranja / Lane Late / Primera
$ws.Cells.Item(1061, 1).Value2 = 8
$ws.Cells.Item(1061, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1061, 3).Value2 = "Coquimbo"
$ws.Cells.Item(1061, 4).Value2 = 45239
$ws.Cells.Item(1061, 5).Value2 = 4
$ws.Cells.Item(1061, 6).Value2 = "Fruta"
$ws.Cells.Item(1061, 7).Value2 = 100102
$ws.Cells.Item(1061, 8).Value2 = "Cítricos"
$ws.Cells.Item(1061, 9).Value2 = 100102005
$ws.Cells.Item(1061, 10).Value2 = "Naranja"
$ws.Cells.Item(1061, 11).Value2 = "Lane Late"
$ws.Cells.Item(1061, 12).Value2 = "Primera"
$ws.Cells.Item(1061, 13).Value2 = 10
$ws.Cells.Item(1061, 14).Value2 = 180000
$ws.Cells.Item(1061, 15).Value2 = 190000
$ws.Cells.Item(1061, 16).Value2 = 185000
$ws.Cells.Item(1061, 17).Value2 = "`$/bins (400 kilos)"
$ws.Cells.Item(1061, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(1061, 19).Value2 = 462
$ws.Cells.Item(1061, 20).Value2 = 400
